$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = [double]"1.982372174814036E-15"
$ws.Range("E3").Value = [double]"1.982372174814036E-15"

$ws.Range("D4").Value = [double]"0.0003687707289405372"
$ws.Range("E4").Value = [double]"0.0003687707289405372"

$ws.Range("D5").Value = [double]"5.713336679434647E-24"
$ws.Range("E5").Value = [double]"5.713336679434647E-24"

$ws.Range("D6").Value = [double]"4.259248333022048E-67"
$ws.Range("E6").Value = [double]"4.259248333022048E-67"

$ws.Range("D8").Value = [double]"0.966073947927751"
$ws.Range("E8").Value = [double]"0.03392605207224897"

$ws.Range("D9").Value = [double]"0.9999999908846561"
$ws.Range("E9").Value = [double]"9.115343924115393E-09"

$ws.Range("D11").Value = [double]"0.9999999999999991"
$ws.Range("E11").Value = [double]"8.881784197001252E-16"
$ws.Range("F11").Value = [double]"4.699339389801025"
